# Apply hybrid bold + color highlighting to quantitative metrics
# (percentages, dollar amounts, large numbers) in specific bullet points
# of the resume, per the target diff.

$d = $word.ActiveDocument

# Color used for all highlighted metrics: 2C3E50 (stored as BGR int for Font.Color)
$metricColor = 0x503E2C

function Highlight-Metric($paragraph, [string]$text) {
    $rng = $paragraph.Range
    $rng.Find.Execute($text) | Out-Null
    $rng.Font.Bold = 1
    $rng.Font.Color = $metricColor
}

# Paragraph 9: Siege Analytics bullet 1 - "23%" and "64%"
$p9 = $d.Paragraphs.Item(9)
Highlight-Metric $p9 "23%"
Highlight-Metric $p9 "64%"

# Paragraph 11: Siege Analytics bullet 3 - "87%", "71%", "±4.2%", "±2.1%"
$p11 = $d.Paragraphs.Item(11)
Highlight-Metric $p11 "87%"
Highlight-Metric $p11 "71%"
Highlight-Metric $p11 "±4.2%"
Highlight-Metric $p11 "±2.1%"

# Paragraph 31: Senior Analyst bullet 3 - "1,200"
$p31 = $d.Paragraphs.Item(31)
Highlight-Metric $p31 "1,200"

# Paragraph 46: Programmer bullet 3 - "$400M" and "$1B"
$p46 = $d.Paragraphs.Item(46)
Highlight-Metric $p46 "$400M"
Highlight-Metric $p46 "$1B"

# Paragraph 63: KEY ACHIEVEMENTS bullet 2 - "73.5%" and "$4.7M"
$p63 = $d.Paragraphs.Item(63)
Highlight-Metric $p63 "73.5%"
Highlight-Metric $p63 "$4.7M"

# Paragraph 65: KEY ACHIEVEMENTS bullet 4 - "87%" and "71%"
$p65 = $d.Paragraphs.Item(65)
Highlight-Metric $p65 "87%"
Highlight-Metric $p65 "71%"

Write-Output "Done applying metric highlighting"
